$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Edificio Mario de Freitas entry ---
# Leading apostrophe preserves the quotePrefix-as-text styling already on
# column B (matches how the original numeric-looking IDs were entered).
$ws.Range("B2").Value = "'0620"
$ws.Range("C2").Value = "Edificio Mário de Freitas"
$ws.Range("D2").Value = "Cliente cobrando troca de bateria da central."
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "Pendente"

# --- Row 3: Escola Odilon Santiago entry ---
$ws.Range("B3").Value = "'0224"
$ws.Range("C3").Value = "Escola Odilon Santiago"
$ws.Range("D3").Value = "Querem refazer relação de usuários mas não querem refazer remoto."
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = "Pendente"

# --- Row 4: Moda Mil entry ---
$ws.Range("B4").Value = "'0236"
$ws.Range("C4").Value = "Moda Mil"
$ws.Range("D4").Value = "Central do cliente ficando sem comunicação e câmeras segue on, cliente pedindo reparo."
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G4").Value = "Pendente"

# --- Row 5: Mf Eventos entry ---
$ws.Range("B5").Value = "'0134"
$ws.Range("C5").Value = "Mf Eventos"
$ws.Range("D5").Value = "Ambos os dvr sem imagem (novamente), irei pedir o Giovani pra ver o que pode tá acontecnedo."
$ws.Range("E5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = "Pendente"

# --- Row 2 Ordens Abertas column, updated last ---
$ws.Range("H2").Value = "Maxvel: 36 / Forte: 19"

# --- Rows 6-12 are wiped entirely and their row heights reset to default ---
$ws.Range("A6:I12").ClearContents()
for ($r = 6; $r -le 12; $r++) {
    $ws.Rows.Item($r).AutoFit()
}

# --- Update the active selection to H2 ---
[void]$ws.Range("H2").Select()
